$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new dataset values
$row2 = @(45150.50694444445, 18.737, 12.527, 4.029, 39.827, 31.89, 14.745, 46.499, 22.687, 9.44, 14.208, 15.679, 16.314, 4.706, 14.663, 20.445, 12.568, 3.441, 2.244, 215.823, 40.754, 13.534, 26.795, 13.632, 3.069, 23.634, 11.955, 10.86, 12.744, 16.181, 3.455, 41.233, 7.422, 16.92)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @(45150.51388888889, 15.374, 10.833, 1.775, 33.225, 26.821, 12.099, 46.928, 18.615, 8.027, 11.826, 13.288, 13.877, 3.865, 12.031, 16.943, 10.387, 1.54, 1.039, 175.806, 33.682, 11.105, 22.313, 11.573, 2.181, 22.799, 9.809, 8.883, 10.415, 13.793, 1.265, 42.715, 6.128, 13.884)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @(45150.52083333334, 9.128, 6.393, 1.077, 19.758, 15.852, 7.184, 30.243, 11.053, 4.756, 6.92, 7.923, 8.258, 2.297, 7.143, 10.049, 6.248, 0.986, 0.606, 101.413, 20.074, 6.594, 13.24, 6.866, 1.349, 14.296, 5.824, 5.309, 6.217, 8.221, 0.773, 27.563, 3.604, 8.244)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @(45150.52777777778, 14.89, 10.84, 1.05, 32.34, 26.35, 11.72, 44.19, 18.03, 7.93, 11.7, 12.98, 13.64, 3.74, 11.65, 16.52, 9.92, 0.8, 0.67, 170.07, 32.55, 10.76, 21.77, 11.4, 1.84, 21.54, 9.5, 8.49, 9.97, 13.58, 0.56, 39.96, 6.01, 13.45)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# Remove the now-unused 6th row (dataset shrank by one row)
$ws.Rows.Item(6).Delete()

# Widen a set of columns from 7 to 8 characters to match the new layout
$wideColumns = @(3, 6, 7, 11, 15, 17, 22, 23, 24, 27, 29)
foreach ($col in $wideColumns) {
    $ws.Columns.Item($col).ColumnWidth = 7.17
}
